# Auto-generated edit script applying the Cerberus_Profits market-data refresh
# Updates cell values in columns H-N across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (120 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 596.3
$ws.Range("I4").Value = 629.625
$ws.Range("J4").Value = 463
$ws.Range("K4").Value = 629.625
$ws.Range("L4").Value = 463
$ws.Range("M4").Value = -515.625
$ws.Range("N4").Value = -691
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null
$ws.Range("H33").Value = 1172.3846
$ws.Range("I33").Value = 1610.7778
$ws.Range("J33").Value = 186
$ws.Range("K33").Value = 1610.7778
$ws.Range("L33").Value = 186
$ws.Range("M33").Value = -1381.7778
$ws.Range("N33").Value = -644
$ws.Range("H34").Value = 14125
$ws.Range("I34").Value = 14125
$ws.Range("K34").Value = 14125
$ws.Range("M34").Value = -13922
$ws.Range("H36").Value = 14125
$ws.Range("I36").Value = 14125
$ws.Range("K36").Value = 14125
$ws.Range("M36").Value = -13410
$ws.Range("H38").Value = 1452.5555
$ws.Range("I38").Value = 399.66666
$ws.Range("J38").Value = 1979
$ws.Range("K38").Value = 1198.99998
$ws.Range("L38").Value = 5937
$ws.Range("M38").Value = -826.99998
$ws.Range("N38").Value = -6681
$ws.Range("H39").Value = 1368.3334
$ws.Range("J39").Value = 9499
$ws.Range("L39").Value = 28497
$ws.Range("N39").Value = -29089
$ws.Range("H40").Value = 2810.5
$ws.Range("I40").Value = 2369
$ws.Range("K40").Value = 2369
$ws.Range("M40").Value = -2194
$ws.Range("H42").Value = 657
$ws.Range("I42").Value = 674.6667
$ws.Range("K42").Value = 2024.0001
$ws.Range("M42").Value = -1794.0001
$ws.Range("H43").Value = 2382
$ws.Range("J43").Value = 2382
$ws.Range("L43").Value = 2382
$ws.Range("N43").Value = -2520
$ws.Range("H74").Value = 5533.878
$ws.Range("I74").Value = 4516.731
$ws.Range("K74").Value = 4516.731
$ws.Range("M74").Value = -3580.731
$ws.Range("H77").Value = 5533.878
$ws.Range("I77").Value = 4516.731
$ws.Range("K77").Value = 22583.655
$ws.Range("M77").Value = -17903.655
$ws.Range("H92").Value = 3204.4167
$ws.Range("I92").Value = 2424.2856
$ws.Range("K92").Value = 2424.2856
$ws.Range("M92").Value = -1176.2856
$ws.Range("H98").Value = 2441.8438
$ws.Range("I98").Value = 1990.5358
$ws.Range("J98").Value = 5601
$ws.Range("K98").Value = 1990.5358
$ws.Range("L98").Value = 5601
$ws.Range("M98").Value = -492.5358000000001
$ws.Range("N98").Value = -8597
$ws.Range("H103").Value = 1576.7273
$ws.Range("I103").Value = 2431.8
$ws.Range("K103").Value = 7295.400000000001
$ws.Range("M103").Value = -6709.400000000001
$ws.Range("H113").Value = 6280.8
$ws.Range("I113").Value = 6279.7646
$ws.Range("J113").Value = 6283
$ws.Range("K113").Value = 6279.7646
$ws.Range("L113").Value = 6283
$ws.Range("M113").Value = -3025.7646
$ws.Range("N113").Value = -12791
$ws.Range("H116").Value = 13605.667
$ws.Range("J116").Value = 12928
$ws.Range("L116").Value = 12928
$ws.Range("N116").Value = -19812
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = $null
$ws.Range("H122").Value = 2441.8438
$ws.Range("I122").Value = 1990.5358
$ws.Range("J122").Value = 5601
$ws.Range("K122").Value = 5971.607400000001
$ws.Range("L122").Value = 16803
$ws.Range("M122").Value = -3521.607400000001
$ws.Range("N122").Value = -21703
$ws.Range("H132").Value = 3232.375
$ws.Range("I132").Value = 2980.5518
$ws.Range("K132").Value = 8941.6554
$ws.Range("M132").Value = -6411.6554
$ws.Range("H137").Value = 2218.6843
$ws.Range("I137").Value = 1891.6154
$ws.Range("K137").Value = 5674.8462
$ws.Range("M137").Value = -3124.8462
$ws.Range("H138").Value = 4307.4517
$ws.Range("I138").Value = 5027.077
$ws.Range("J138").Value = 3787.7222
$ws.Range("K138").Value = 15081.231
$ws.Range("L138").Value = 11363.1666
$ws.Range("M138").Value = -9941.231
$ws.Range("N138").Value = -21643.1666
$ws.Range("H141").Value = 5107.76
$ws.Range("I141").Value = 2999.8
$ws.Range("J141").Value = 13539.6
$ws.Range("K141").Value = 8999.400000000001
$ws.Range("L141").Value = 40618.8
$ws.Range("M141").Value = -3819.400000000001
$ws.Range("N141").Value = -50978.8

# --- Sheet: ARM (94 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 22499.5
$ws.Range("J24").Value = 22499.5
$ws.Range("L24").Value = 22499.5
$ws.Range("N24").Value = -23247.5
$ws.Range("H26").Value = 4990
$ws.Range("I26").Value = 4990
$ws.Range("K26").Value = 4990
$ws.Range("M26").Value = -4660
$ws.Range("H45").Value = 3279
$ws.Range("I45").Value = 1199
$ws.Range("K45").Value = 1199
$ws.Range("M45").Value = -822
$ws.Range("H46").Value = 49664.332
$ws.Range("I46").Value = 49664.332
$ws.Range("K46").Value = 49664.332
$ws.Range("M46").Value = -49345.332
$ws.Range("H61").Value = 9970.632
$ws.Range("I61").Value = 5461.8667
$ws.Range("K61").Value = 5461.8667
$ws.Range("M61").Value = -5249.8667
$ws.Range("H63").Value = 1340.0834
$ws.Range("I63").Value = 1371.1818
$ws.Range("J63").Value = 998
$ws.Range("K63").Value = 1371.1818
$ws.Range("L63").Value = 998
$ws.Range("M63").Value = -685.1818000000001
$ws.Range("N63").Value = -2370
$ws.Range("H66").Value = 1340.0834
$ws.Range("I66").Value = 1371.1818
$ws.Range("J66").Value = 998
$ws.Range("K66").Value = 6855.909000000001
$ws.Range("L66").Value = 4990
$ws.Range("M66").Value = -3423.909000000001
$ws.Range("N66").Value = -11854
$ws.Range("H74").Value = 3706.0386
$ws.Range("I74").Value = 1288.3077
$ws.Range("J74").Value = 6123.769
$ws.Range("K74").Value = 1288.3077
$ws.Range("L74").Value = 6123.769
$ws.Range("M74").Value = -414.3077000000001
$ws.Range("N74").Value = -7871.769
$ws.Range("H77").Value = 3706.0386
$ws.Range("I77").Value = 1288.3077
$ws.Range("J77").Value = 6123.769
$ws.Range("K77").Value = 6441.538500000001
$ws.Range("L77").Value = 30618.845
$ws.Range("M77").Value = -2073.538500000001
$ws.Range("N77").Value = -39354.845
$ws.Range("H88").Value = 26348.625
$ws.Range("J88").Value = 27677.2
$ws.Range("L88").Value = 27677.2
$ws.Range("N88").Value = -28489.2
$ws.Range("H91").Value = 26348.625
$ws.Range("J91").Value = 27677.2
$ws.Range("L91").Value = 27677.2
$ws.Range("N91").Value = -30485.2
$ws.Range("H97").Value = 1716.0625
$ws.Range("I97").Value = 663.2727
$ws.Range("K97").Value = 663.2727
$ws.Range("M97").Value = -167.2727
$ws.Range("H100").Value = 22499.5
$ws.Range("J100").Value = 22499.5
$ws.Range("L100").Value = 22499.5
$ws.Range("N100").Value = -24663.5
$ws.Range("H102").Value = 1451.1538
$ws.Range("I102").Value = 1313.24
$ws.Range("K102").Value = 1313.24
$ws.Range("M102").Value = 308.76
$ws.Range("H110").Value = 102467.73
$ws.Range("I110").Value = 139989.12
$ws.Range("J110").Value = 2410.6667
$ws.Range("K110").Value = 139989.12
$ws.Range("L110").Value = 2410.6667
$ws.Range("M110").Value = -137944.12
$ws.Range("N110").Value = -6500.6667
$ws.Range("H114").Value = 21859
$ws.Range("J114").Value = 21859
$ws.Range("L114").Value = 21859
$ws.Range("N114").Value = -30537
$ws.Range("H122").Value = 2711.6667
$ws.Range("I122").Value = 2950
$ws.Range("J122").Value = 2592.5
$ws.Range("K122").Value = 8850
$ws.Range("L122").Value = 7777.5
$ws.Range("M122").Value = -6400
$ws.Range("N122").Value = -12677.5
$ws.Range("H132").Value = 2348.1714
$ws.Range("I132").Value = 1861.1482
$ws.Range("K132").Value = 5583.444600000001
$ws.Range("M132").Value = -3053.444600000001
$ws.Range("H136").Value = 9970.632
$ws.Range("I136").Value = 5461.8667
$ws.Range("K136").Value = 16385.6001
$ws.Range("M136").Value = -13835.6001

# --- Sheet: BSM (48 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2718.6428
$ws.Range("I3").Value = 1962.125
$ws.Range("J3").Value = 3727.3333
$ws.Range("K3").Value = 1962.125
$ws.Range("L3").Value = 3727.3333
$ws.Range("M3").Value = -1848.125
$ws.Range("N3").Value = -3955.3333
$ws.Range("H20").Value = 1556
$ws.Range("I20").Value = 1166.3334
$ws.Range("K20").Value = 1166.3334
$ws.Range("M20").Value = -919.3334
$ws.Range("H22").Value = 662.75
$ws.Range("I22").Value = 467.16666
$ws.Range("J22").Value = 1249.5
$ws.Range("K22").Value = 467.16666
$ws.Range("L22").Value = 1249.5
$ws.Range("M22").Value = -294.16666
$ws.Range("N22").Value = -1595.5
$ws.Range("H80").Value = 843.6667
$ws.Range("I80").Value = 975.75
$ws.Range("J80").Value = 738
$ws.Range("K80").Value = 975.75
$ws.Range("L80").Value = 738
$ws.Range("M80").Value = 22.25
$ws.Range("N80").Value = -2734
$ws.Range("H83").Value = 843.6667
$ws.Range("I83").Value = 975.75
$ws.Range("J83").Value = 738
$ws.Range("K83").Value = 4878.75
$ws.Range("L83").Value = 3690
$ws.Range("M83").Value = 113.25
$ws.Range("N83").Value = -13674
$ws.Range("H86").Value = 16629.223
$ws.Range("I86").Value = 2348.5
$ws.Range("K86").Value = 2348.5
$ws.Range("M86").Value = -1225.5
$ws.Range("H89").Value = 16629.223
$ws.Range("I89").Value = 2348.5
$ws.Range("K89").Value = 11742.5
$ws.Range("M89").Value = -6126.5
$ws.Range("H126").Value = 147495
$ws.Range("J126").Value = 147495
$ws.Range("L126").Value = 147495
$ws.Range("N126").Value = -157375
$ws.Range("H134").Value = 9230.389
$ws.Range("I134").Value = 12813.7
$ws.Range("K134").Value = 38441.10000000001
$ws.Range("M134").Value = -35906.10000000001

# --- Sheet: CRP (51 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 39994.75
$ws.Range("J20").Value = 39994.75
$ws.Range("L20").Value = 39994.75
$ws.Range("N20").Value = -40466.75
$ws.Range("H22").Value = 726.4286
$ws.Range("I22").Value = 503.5
$ws.Range("J22").Value = 815.6
$ws.Range("K22").Value = 503.5
$ws.Range("L22").Value = 815.6
$ws.Range("M22").Value = -153.5
$ws.Range("N22").Value = -1515.6
$ws.Range("H30").Value = 39994.75
$ws.Range("J30").Value = 39994.75
$ws.Range("L30").Value = 39994.75
$ws.Range("N30").Value = -40176.75
$ws.Range("H31").Value = 3049.7827
$ws.Range("J31").Value = 3948.0715
$ws.Range("L31").Value = 3948.0715
$ws.Range("N31").Value = -4538.0715
$ws.Range("H34").Value = 3049.7827
$ws.Range("J34").Value = 3948.0715
$ws.Range("L34").Value = 3948.0715
$ws.Range("N34").Value = -4352.0715
$ws.Range("H92").Value = 33316.668
$ws.Range("J92").Value = 33316.668
$ws.Range("L92").Value = 33316.668
$ws.Range("N92").Value = -38308.668
$ws.Range("H99").Value = 2845.389
$ws.Range("I99").Value = 2692.25
$ws.Range("K99").Value = 2692.25
$ws.Range("M99").Value = -1194.25
$ws.Range("H122").Value = 3602.7273
$ws.Range("I122").Value = 3059.7144
$ws.Range("K122").Value = 9179.143199999999
$ws.Range("M122").Value = -6729.143199999999
$ws.Range("H126").Value = 2845.389
$ws.Range("I126").Value = 2692.25
$ws.Range("K126").Value = 8076.75
$ws.Range("M126").Value = -5606.75
$ws.Range("H128").Value = 39994.75
$ws.Range("J128").Value = 39994.75
$ws.Range("L128").Value = 39994.75
$ws.Range("N128").Value = -49954.75
$ws.Range("H132").Value = 2444.577
$ws.Range("I132").Value = 2410.8333
$ws.Range("K132").Value = 7232.499899999999
$ws.Range("M132").Value = -4702.499899999999
$ws.Range("H134").Value = 4545.7427
$ws.Range("J134").Value = 6369.1816
$ws.Range("L134").Value = 19107.5448
$ws.Range("N134").Value = -24177.5448

# --- Sheet: CUL (67 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2037.5
$ws.Range("I2").Value = 4000
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = -23887
$ws.Range("N2").Value = -676
$ws.Range("H38").Value = 539.9474
$ws.Range("I38").Value = 799.6
$ws.Range("J38").Value = 251.44444
$ws.Range("K38").Value = 2398.8
$ws.Range("L38").Value = 754.33332
$ws.Range("M38").Value = -2051.8
$ws.Range("N38").Value = -1448.33332
$ws.Range("H69").Value = 4948.8335
$ws.Range("J69").Value = 3399.5
$ws.Range("L69").Value = 10198.5
$ws.Range("N69").Value = -11820.5
$ws.Range("H72").Value = 4948.8335
$ws.Range("J72").Value = 3399.5
$ws.Range("L72").Value = 30595.5
$ws.Range("N72").Value = -38707.5
$ws.Range("H87").Value = 13331.777
$ws.Range("I87").Value = 9998
$ws.Range("K87").Value = 29994
$ws.Range("M87").Value = -28746
$ws.Range("H90").Value = 13331.777
$ws.Range("I90").Value = 9998
$ws.Range("K90").Value = 89982
$ws.Range("M90").Value = -83742
$ws.Range("H107").Value = 1446.5555
$ws.Range("I107").Value = 389.33334
$ws.Range("J107").Value = 1658
$ws.Range("K107").Value = 1168.00002
$ws.Range("L107").Value = 4974
$ws.Range("M107").Value = 751.99998
$ws.Range("N107").Value = -8814
$ws.Range("H113").Value = 2080
$ws.Range("J113").Value = 2107
$ws.Range("L113").Value = 6321
$ws.Range("N113").Value = -10661
$ws.Range("H118").Value = 100
$ws.Range("I118").Value = 100
$ws.Range("K118").Value = 300
$ws.Range("M118").Value = 943
$ws.Range("H119").Value = 4415.2856
$ws.Range("I119").Value = 4817.8335
$ws.Range("K119").Value = 14453.5005
$ws.Range("M119").Value = -9615.500499999998
$ws.Range("H121").Value = 5333.9
$ws.Range("I121").Value = 3280
$ws.Range("J121").Value = 6214.143
$ws.Range("K121").Value = 9840
$ws.Range("L121").Value = 18642.429
$ws.Range("M121").Value = -8530
$ws.Range("N121").Value = -21262.429
$ws.Range("H128").Value = 194980
$ws.Range("I128").Value = 194980
$ws.Range("K128").Value = 584940
$ws.Range("M128").Value = -579960
$ws.Range("H132").Value = 1393
$ws.Range("I132").Value = 1450.6
$ws.Range("J132").Value = 1105
$ws.Range("K132").Value = 13055.4
$ws.Range("L132").Value = 9945
$ws.Range("M132").Value = -10525.4
$ws.Range("N132").Value = -15005

# --- Sheet: GSM (41 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 30666.666
$ws.Range("I58").Value = 24000
$ws.Range("K58").Value = 24000
$ws.Range("M58").Value = -23723
$ws.Range("H97").Value = 5552.25
$ws.Range("I97").Value = 5136.125
$ws.Range("J97").Value = 6384.5
$ws.Range("K97").Value = 5136.125
$ws.Range("L97").Value = 6384.5
$ws.Range("M97").Value = -4640.125
$ws.Range("N97").Value = -7376.5
$ws.Range("H102").Value = 3785.4546
$ws.Range("I102").Value = 3736.5
$ws.Range("K102").Value = 3736.5
$ws.Range("M102").Value = -2114.5
$ws.Range("H113").Value = 3823.1667
$ws.Range("I113").Value = 3236.6667
$ws.Range("J113").Value = 4409.6665
$ws.Range("K113").Value = 3236.6667
$ws.Range("L113").Value = 4409.6665
$ws.Range("M113").Value = -1066.6667
$ws.Range("N113").Value = -8749.6665
$ws.Range("H122").Value = 3805.4443
$ws.Range("I122").Value = 3646.111
$ws.Range("K122").Value = 10938.333
$ws.Range("M122").Value = -8488.332999999999
$ws.Range("H123").Value = 64166.332
$ws.Range("J123").Value = 64166.332
$ws.Range("L123").Value = 64166.332
$ws.Range("N123").Value = -69066.332
$ws.Range("H126").Value = 4154.6816
$ws.Range("I126").Value = 2898.077
$ws.Range("K126").Value = 8694.231
$ws.Range("M126").Value = -6224.231
$ws.Range("H132").Value = 4680.871
$ws.Range("I132").Value = 4922.8076
$ws.Range("J132").Value = 3422.8
$ws.Range("K132").Value = 14768.4228
$ws.Range("L132").Value = 10268.4
$ws.Range("M132").Value = -12238.4228
$ws.Range("N132").Value = -15328.4

# --- Sheet: LTW (86 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1900.75
$ws.Range("I12").Value = 1900.75
$ws.Range("K12").Value = 1900.75
$ws.Range("M12").Value = -1730.75
$ws.Range("H22").Value = 1489.6
$ws.Range("I22").Value = 483.33334
$ws.Range("K22").Value = 483.33334
$ws.Range("M22").Value = -188.33334
$ws.Range("H27").Value = 1489.6
$ws.Range("I27").Value = 483.33334
$ws.Range("K27").Value = 483.33334
$ws.Range("M27").Value = -376.33334
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -664
$ws.Range("N35").Value = $null
$ws.Range("H40").Value = 2031.2941
$ws.Range("I40").Value = 1643.6666
$ws.Range("K40").Value = 1643.6666
$ws.Range("M40").Value = -1507.6666
$ws.Range("H46").Value = 9092654
$ws.Range("I46").Value = 16667417
$ws.Range("J46").Value = 2939.2
$ws.Range("K46").Value = 16667417
$ws.Range("L46").Value = 2939.2
$ws.Range("M46").Value = -16667229
$ws.Range("N46").Value = -3315.2
$ws.Range("H61").Value = 2292.25
$ws.Range("I61").Value = 2045.6364
$ws.Range("K61").Value = 2045.6364
$ws.Range("M61").Value = -1843.6364
$ws.Range("H68").Value = 2172.5
$ws.Range("I68").Value = 2172.5
$ws.Range("K68").Value = 2172.5
$ws.Range("M68").Value = -1423.5
$ws.Range("H71").Value = 2172.5
$ws.Range("I71").Value = 2172.5
$ws.Range("K71").Value = 10862.5
$ws.Range("M71").Value = -7118.5
$ws.Range("H82").Value = 1741.1111
$ws.Range("I82").Value = 1399.25
$ws.Range("J82").Value = 2424.8333
$ws.Range("K82").Value = 1399.25
$ws.Range("L82").Value = 2424.8333
$ws.Range("M82").Value = -1038.25
$ws.Range("N82").Value = -3146.8333
$ws.Range("H85").Value = 1741.1111
$ws.Range("I85").Value = 1399.25
$ws.Range("J85").Value = 2424.8333
$ws.Range("K85").Value = 1399.25
$ws.Range("L85").Value = 2424.8333
$ws.Range("M85").Value = -151.25
$ws.Range("N85").Value = -4920.8333
$ws.Range("H93").Value = 1911.1538
$ws.Range("I93").Value = 1820.4166
$ws.Range("K93").Value = 1820.4166
$ws.Range("M93").Value = -572.4166
$ws.Range("H100").Value = 860.75
$ws.Range("I100").Value = 698
$ws.Range("K100").Value = 698
$ws.Range("M100").Value = -157
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954
$ws.Range("H113").Value = 2292.25
$ws.Range("I113").Value = 2045.6364
$ws.Range("K113").Value = 2045.6364
$ws.Range("M113").Value = 124.3635999999999
$ws.Range("H132").Value = 3130.2334
$ws.Range("I132").Value = 2571.7917
$ws.Range("J132").Value = 5364
$ws.Range("K132").Value = 7715.375100000001
$ws.Range("L132").Value = 16092
$ws.Range("M132").Value = -5185.375100000001
$ws.Range("N132").Value = -21152
$ws.Range("H136").Value = 1945.305
$ws.Range("I136").Value = 1235.875
$ws.Range("J136").Value = 3438.842
$ws.Range("K136").Value = 3707.625
$ws.Range("L136").Value = 10316.526
$ws.Range("M136").Value = -1157.625
$ws.Range("N136").Value = -15416.526

# --- Sheet: WVR (53 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 4000
$ws.Range("I9").Value = 4000
$ws.Range("K9").Value = 4000
$ws.Range("M9").Value = -3860
$ws.Range("H62").Value = 15784
$ws.Range("I62").Value = 6513.6665
$ws.Range("J62").Value = 19757
$ws.Range("K62").Value = 6513.6665
$ws.Range("L62").Value = 19757
$ws.Range("M62").Value = -5889.6665
$ws.Range("N62").Value = -21005
$ws.Range("H65").Value = 15784
$ws.Range("I65").Value = 6513.6665
$ws.Range("J65").Value = 19757
$ws.Range("K65").Value = 32568.3325
$ws.Range("L65").Value = 98785
$ws.Range("M65").Value = -29448.3325
$ws.Range("N65").Value = -105025
$ws.Range("H74").Value = 25933
$ws.Range("J74").Value = 25933
$ws.Range("L74").Value = 25933
$ws.Range("N74").Value = -27805
$ws.Range("H76").Value = 70000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100630
$ws.Range("H77").Value = 25933
$ws.Range("J77").Value = 25933
$ws.Range("L77").Value = 77799
$ws.Range("N77").Value = -87159
$ws.Range("H79").Value = 70000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102184
$ws.Range("H92").Value = 34775
$ws.Range("J92").Value = 34775
$ws.Range("L92").Value = 34775
$ws.Range("N92").Value = -39767
$ws.Range("H96").Value = 3443.3333
$ws.Range("I96").Value = 1999
$ws.Range("K96").Value = 1999
$ws.Range("M96").Value = -626
$ws.Range("H117").Value = 22469.666
$ws.Range("J117").Value = 22469.666
$ws.Range("L117").Value = 22469.666
$ws.Range("N117").Value = -31647.666
$ws.Range("H132").Value = 1535.6364
$ws.Range("I132").Value = 1489.8108
$ws.Range("J132").Value = 1777.8572
$ws.Range("K132").Value = 4469.4324
$ws.Range("L132").Value = 5333.571599999999
$ws.Range("M132").Value = -1939.4324
$ws.Range("N132").Value = -10393.5716

Write-Host "Applied 560 cell updates."